$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that were removed entirely in the diff (C2, E2, C3)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

# Apply tiny floating point precision corrections (recomputed forecast values)
$ws.Range("E3").Value = 0.5079568386449518
$ws.Range("C4").Value = -0.9140166223623458
$ws.Range("E4").Value = 5.639535270494123
$ws.Range("E6").Value = -0.6955733540840225
$ws.Range("C8").Value = -1.479696720105184
$ws.Range("E9").Value = -2.092856741436244
$ws.Range("E10").Value = -1.259568900987029
$ws.Range("C11").Value = 2.192778679161966
$ws.Range("C12").Value = 3.408364488606752
$ws.Range("E12").Value = -0.3010260522302244
$ws.Range("E13").Value = 3.825329033908798
$ws.Range("C15").Value = 1.666553973046025
$ws.Range("E15").Value = -7.134843267358049
$ws.Range("C16").Value = 1.879266440112781
$ws.Range("E16").Value = -6.539839435602913
$ws.Range("C17").Value = -2.620683231370935
$ws.Range("E17").Value = -5.419975784955122
$ws.Range("C18").Value = -3.036556262700263
